# Update the "kbase_text" worksheet:
#  - column B header changes from "类型" to "标题"
#  - the per-row "TEXT" markers become distinct "标题1" / "标题2" labels
#  - the active selection moves from C8 to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "标题"
$ws.Range("B2").Value = "标题1"
$ws.Range("B3").Value = "标题2"

$ws.Range("B4").Select()
